$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the data columns (D = Price, E = Volume(1h)) as Text so the
# new values (e.g. "300.82", "-0.22%", "2,116.77%") are stored as literal
# strings -- matching the original inline-string cells -- instead of being
# auto-converted to numbers/percentages by the Excel input parser.
$numRng = $ws.Range("D2:E50")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = '300.82'
$ws.Range("E2").Value = '-0.22%'
$ws.Range("D3").Value = '32.48'
$ws.Range("E3").Value = '2.10%'
$ws.Range("D4").Value = '5.062'
$ws.Range("E4").Value = '-1.32%'
$ws.Range("D5").Value = '0.07677'
$ws.Range("E5").Value = '-2.10%'
$ws.Range("D6").Value = '2.125'
$ws.Range("E6").Value = '-5.84%'
$ws.Range("D7").Value = '7.852'
$ws.Range("E7").Value = '0.54%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9208'
$ws.Range("E8").Value = '-0.80%'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '0.1760'
$ws.Range("E9").Value = '-0.82%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.07887'
$ws.Range("E10").Value = '2.70%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.08469'
$ws.Range("E11").Value = '-4.89%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03059'
$ws.Range("E12").Value = '-1.36%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09993'
$ws.Range("E13").Value = '-0.27%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001525'
$ws.Range("E14").Value = '1.19%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005817'
$ws.Range("E15").Value = '0.23%'
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = '0.007498'
$ws.Range("E16").Value = '2,116.77%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.461'
$ws.Range("E17").Value = '-0.59%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '3.784'
$ws.Range("E18").Value = '-0.63%'
$ws.Range("E19").Value = '-4.44%'
$ws.Range("D20").Value = '0.3340'
$ws.Range("E20").Value = '1.47%'
$ws.Range("D21").Value = '0.1317'
$ws.Range("E21").Value = '-0.81%'
$ws.Range("D22").Value = '4.269'
$ws.Range("E22").Value = '-1.22%'
$ws.Range("D24").Value = '0.04526'
$ws.Range("E24").Value = '-1.61%'
$ws.Range("D25").Value = '0.001236'
$ws.Range("E25").Value = '-1.18%'
$ws.Range("D26").Value = '0.004831'
$ws.Range("E26").Value = '7.69%'
$ws.Range("D27").Value = '0.0001252'
$ws.Range("E27").Value = '0.24%'
$ws.Range("D39").Value = '0.01708'
$ws.Range("E39").Value = '-4.21%'
$ws.Range("D40").Value = '0.04675'
$ws.Range("E40").Value = '-2.48%'
$ws.Range("D41").Value = '0.007458'
$ws.Range("E41").Value = '1.41%'
$ws.Range("D42").Value = '0.1354'
$ws.Range("E42").Value = '-0.84%'
$ws.Range("D43").Value = '0.002334'
$ws.Range("E43").Value = '6.65%'
$ws.Range("E44").Value = '7.04%'
$ws.Range("D45").Value = '0.00006210'
$ws.Range("E45").Value = '-0.87%'
$ws.Range("E46").Value = '0.07%'
$ws.Range("B47").Value = 'BOLO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D47").Value = '1.046'
$ws.Range("E47").Value = '-2.56%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '0.003000'
$ws.Range("E48").Value = '-62.44%'
$ws.Range("D49").Value = '0.00002100'
$ws.Range("E49").Value = '0.07%'
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").Value = '0.07%'

# Restore the default (unstyled) cell formatting now that the text values
# are committed, so no stray number-format style is left behind.
$numRng.Style = "Normal"
